# Rotate the comma-separated "Recorded By" list in column G so that the
# last entry moves to the front (e.g. "System, a@b.com" -> "a@b.com, System").
# Only rows whose value actually contains a comma (i.e. more than one
# recorder) are touched; single-value cells are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($null -ne $val -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ", "
        if ($parts.Length -gt 1) {
            $last = $parts[$parts.Length - 1]
            $rest = $parts[0..($parts.Length - 2)]
            $newParts = @($last) + $rest
            $newVal = [string]::Join(", ", $newParts)
            $cell.Value2 = $newVal
        }
    }
}
